$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new text value, as scraped from the source diff.
# These cells hold plain text (prices / percentage strings), not true
# numbers, so each value is forced to Text format before being written
# and the cell's style is reset back to Normal afterwards so that no
# stray formatting is introduced.
$updates = @{
    "D2" = "305.71"
    "E2" = "-0.72%"
    "D3" = "38.89"
    "E3" = "6.99%"
    "D4" = "5.108"
    "E4" = "0.94%"
    "E5" = "-0.26%"
    "D6" = "1.938"
    "E6" = "-3.73%"
    "D7" = "4.202"
    "E7" = "0.92%"
    "D8" = "8.037"
    "E8" = "2.64%"
    "D9" = "0.9256"
    "E9" = "-0.20%"
    "D10" = "0.1411"
    "E10" = "-5.47%"
    "D11" = "0.1914"
    "E11" = "-1.06%"
    "D12" = "0.09024"
    "E12" = "-0.95%"
    "D13" = "0.03506"
    "E13" = "-0.67%"
    "D14" = "0.09774"
    "E14" = "-0.89%"
    "D15" = "0.001394"
    "E15" = "-1.37%"
    "D16" = "0.005934"
    "E16" = "-5.56%"
    "D17" = "3.761"
    "E17" = "-2.28%"
    "E18" = "-2.12%"
    "D19" = "0.3462"
    "E19" = "0.44%"
    "D20" = "0.1293"
    "E20" = "-2.54%"
    "D21" = "4.673"
    "E21" = "-3.22%"
    "D22" = "0.2416"
    "E22" = "2.98%"
    "D23" = "0.04381"
    "E23" = "-0.32%"
    "D24" = "0.001206"
    "E24" = "-2.55%"
    "E25" = "2.82%"
    "E26" = "-0.19%"
    "D39" = "0.02036"
    "E39" = "-1.18%"
    "D40" = "0.05046"
    "E40" = "-1.05%"
    "E41" = "0.73%"
    "D42" = "0.009724"
    "E42" = "-4.00%"
    "D43" = "0.1342"
    "E43" = "-1.74%"
    "E44" = "-1.60%"
    "D45" = "0.009784"
    "E45" = "1.22%"
    "D46" = "0.00006244"
    "E46" = "-1.04%"
    "E47" = "-0.18%"
    "D48" = "0.002871"
    "E49" = "12.35%"
    "E50" = "-0.18%"
    "E51" = "-0.18%"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
    $cell.Style = "Normal"
}
